# "Rerun failed cases implemented"
# Each of the per-brand sheets records form submissions where column A holds a
# TimeStamp (as text) that is identical for every row of that sheet (one
# submission "run" = one timestamp, fanned out to every row it produced).
# This rerun refreshes the TimeStamp for every existing row in every sheet to
# a new value, without touching the rest of the data (Bike Name / Price /
# Launch Date stay the same).

$wb = $excel.ActiveWorkbook

$newTimestamps = @{
    "Royal Enfield"  = "Feb 16, 2026 12:26:50"
    "Yamaha"         = "Feb 16, 2026 12:27:11"
    "TVS"            = "Feb 16, 2026 12:27:28"
    "Bajaj"          = "Feb 16, 2026 12:27:43"
    "Hero Moto Corp" = "Feb 16, 2026 12:27:59"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($newTimestamps.ContainsKey($name)) {
        $stamp = $newTimestamps[$name]
        $used = $ws.UsedRange
        $lastRow = $used.Rows.Count
        for ($r = 2; $r -le $lastRow; $r++) {
            $ws.Cells.Item($r, 1).Value = $stamp
        }
    }
}
